$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.454.10'
$ws.Range('D3').Value = '1.870.76'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.7060'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3152'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07871'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.66'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('D12').Value = '1.890.26'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.214'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '94.14'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.489'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '29.518.92'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008370'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '256.70'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.51%  '
$ws.Range('D20').Value = '2.141.88'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.624'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1554'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.056'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '161.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.500'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.336'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.253'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.208'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.902'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7480'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.74%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.170'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.718'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01876'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').Value = '1.264.76'
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.749'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8997'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '109.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.89%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '71.78'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.76%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.946'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -8.92%  '
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('D47').Value = '2.038.12'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.817'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.5192'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.509'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06096'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.29%  '
